# Bug Fixed - UI Update
# Fill in the actual OCT figures (previously zeroed placeholders) and
# propagate the corrected running balance (SOLD) through the OCT, NOV
# and DEC rows of the budget table.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Row 12 = OCT, Row 13 = NOV, Row 14 = DEC (rows 1-2 are headers,
# rows 3-11 are JAN..SEPT).
$octRow = 12
$novRow = 13
$decRow = 14

# OCT row: NOMBRES / MONTANT CONSOMME for the two payment buckets, then
# the NOMBRES TOTAL / MONTANT CONSOMME TOTAL columns.
$table.Cell($octRow, 3).Range.Text  = "4173"
$table.Cell($octRow, 4).Range.Text  = "3 755 700,00"
$table.Cell($octRow, 5).Range.Text  = "518"
$table.Cell($octRow, 6).Range.Text  = "644 200,00"
$table.Cell($octRow, 7).Range.Text  = "4691"
$table.Cell($octRow, 8).Range.Text  = "4 399 900,00"

# SOLD (running balance) column updates for OCT, NOV, DEC rows.
$table.Cell($octRow, 9).Range.Text  = "15 379 500,00"
$table.Cell($novRow, 9).Range.Text  = "15 379 500,00"
$table.Cell($decRow, 9).Range.Text  = "15 379 500,00"
